# Update the "as of" date in the confidential disclosure note (shared string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect for the duration of the edit and
# restore protection afterwards.
$ws.Unprotect()

$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Updated Weight (D) and Percent Change (E) figures for each holding row.
$updates = @(
    @(2, 0.06568501266571172, -0.03538554398672089),
    @(3, 0.0399594800747135, -0.02203461401037643),
    @(4, 0.03269052652386122, -0.01615977130151691),
    @(5, 0.03019194542490017, -0.002493443961996467),
    @(6, 0.02764758158844695, -0.01547108933540464),
    @(7, 0.02412797464536265, 0.01382368283776714),
    @(8, 0.1847278392661226, 0.003917727717923647),
    @(9, 0.02404289899141456, 0.01549543005871312),
    @(10, 0.02230138501552005, -0.003407407407407526),
    @(11, 0.02227400977042164, -0.006604226705091198),
    @(12, 0.02043027128014383, -0.006792086679963338),
    @(13, 0.01946285529100497, 0.01084812623274156),
    @(14, 0.01670483867958825, -0.008989642368575357),
    @(15, 0.01680462459455188, -0.02376137512639043),
    @(16, 0.01511491119020136, -0.01240736114580743),
    @(17, 0.01453373788334188, -0.009850830284266809),
    @(18, 0.01433562134661671, -0.0125523012552301),
    @(19, 0.01395657646585027, -0.01308202616405219),
    @(20, 0.01295572806531249, 0.006290377422645266),
    @(21, 0.01226576175853603, 0.01650793650793658),
    @(22, 0.01308080461619319, -0.006801534705266787),
    @(23, 0.01193812412682654, -0.009198734844491185),
    @(24, 0.0127719677993645, 0.01340231584134033),
    @(25, 0.01169123560170763, 0.007569564295879072),
    @(26, 0.009363277797627235, -0.007174782404140156),
    @(27, 0.009710974876176084, -0.01774026310673327),
    @(28, 0.0102200915034087, 0.004017857142857073),
    @(29, 0.01034131349391634, -0.02838484417431775),
    @(30, 0.01007062895833114, -0.01382596469301678),
    @(31, 0.008783441787223713, -0.006287111421585778),
    @(32, 0.01023511642241386, 0.00106063284426372),
    @(33, 0.009340307764383738, 0.01174033149171261),
    @(34, 0.008999847819596538, -0.006240822320117356),
    @(35, 0.009186440007451259, 0.01036136324713133),
    @(36, 0.008462490637794006, -0.01664854011545203),
    @(37, 0.008517005134498611, 0.002147409254640964),
    @(38, 0.008081597141338458, -0.0164987589429112),
    @(39, 0.008772979409068285, 0.003012804418779824),
    @(40, 0.007886391191017424, 0.005496067468966137),
    @(41, 0.007082007372759222, -0.01244057404363075),
    @(42, 0.007567603315266472, -0.02408498872153086),
    @(43, 0.007908063260053671, 0.01551791025475246),
    @(44, 0.007325316663246008, -0.01112531008043283),
    @(45, 0.007353242559826288, -0.004493131926911809),
    @(46, 0.007640525304364839, 0.02100321225599222),
    @(47, 0.007172314215785032, 0.0001316135825215614),
    @(48, 0.007150288156510445, 0.001265182186234837),
    @(49, 0.006827291730148097, -0.01391289319045963),
    @(50, 0.006722589284096396, 0.04414982623246244),
    @(51, 0.006282107430853346, 0.02296533286584568),
    @(52, 0.006494698235101823, 0.01309319057193381),
    @(53, 0.005495108466522586, -0.04101352802233194),
    @(54, 0.005850121343330495, 0.01226333907056798),
    @(55, 0.005949985922791528, 0.008031730292513473),
    @(56, 0.005644881736364914, -0.0009260168940599112),
    @(57, 0.006677042540096445, 0.009943449575871899),
    @(58, 0.005501480290812736, -0.002059025394646641),
    @(59, 0.005076927998295058, 0.008645935016036832),
    @(60, 0.004892931738854842, 0.004758842443729927),
    @(61, 0.004667046634544204, 0.0133494020580327),
    @(62, 0.004654971634191885, 0.02365863962822146),
    @(63, 0.004239701752368498, -0.0001484340210776702),
    @(64, 0.003994897836430937, 0.007403906742281086),
    @(65, 0.003982311516845458, -0.006321112515802807),
    @(66, 0.00370509782797529, -0.009171974522293014),
    @(67, 0.00379005548517727, 0.0004981320049814997),
    @(68, 0.00359359090289769, -0.0008865539320310756),
    @(69, 0.003544622253260436, 0.02751886373723922),
    @(70, 0.003127228430006999, -0.02767017155506368),
    @(71, 0.00305717769506407, -0.03446678760276345),
    @(72, 0.002334250963873137, -0.03272279980453952),
    @(73, 0.002002444113800955, -0.01164777749405832),
    @(74, 0.001980221393282844, -0.0253446152623843),
    @(75, 0.001537969588848087, -0.04475474400286428),
    @(76, 0.001506267796392162, 0.009452684353457252),
    @(77, $null, -0.004375197415930954)
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}

$ws.Protect()
